$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 158, shifting existing rows 158:241 down to 159:242
$ws.Rows(158).Insert()

# Fill in the new row 158 with the new record's data
$ws.Range("A158").Value = 7
$ws.Range("B158").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C158").Value = "Ñuble"
$ws.Range("D158").Value = 44830
$ws.Range("E158").Value = 16
$ws.Range("F158").Value = 100112017
$ws.Range("G158").Value = "Apio"
$ws.Range("H158").Value = "Americana (o)"
$ws.Range("I158").Value = "Primera"
$ws.Range("J158").Value = 120
$ws.Range("K158").Value = 9500
$ws.Range("L158").Value = 10000
$ws.Range("M158").Value = 9750
$ws.Range("N158").Value = "$/docena de matas"
$ws.Range("O158").Value = "Provincia del Elquí"
$ws.Range("P158").Value = 1625
$ws.Range("Q158").Value = 6
$ws.Range("R158").Value = "Hortaliza"
